$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 25.35940266666667
$ws.Range("H2").Value = 76.078208
$ws.Range("I2").Value = 0.005186643687654987
$ws.Range("J2").Value = 0.005186643687654986
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 3694.883852835186
$ws.Range("R2").Value = 33253.95467551668
$ws.Range("S2").Value = 0.00148645979407986
$ws.Range("T2").Value = 0.00148645979407986
$ws.Range("G3").Value = 25.35940266666667
$ws.Range("H3").Value = 76.078208
$ws.Range("I3").Value = 0.005186643687654987
$ws.Range("J3").Value = 0.005186643687654986
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 4280.659629937608
$ws.Range("R3").Value = 38525.93666943847
$ws.Range("S3").Value = 0.001722118660688211
$ws.Range("T3").Value = 0.001722118660688211
$ws.Range("G4").Value = 25.35940266666667
$ws.Range("H4").Value = 76.078208
$ws.Range("I4").Value = 0.005186643687654987
$ws.Range("J4").Value = 0.005186643687654986
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 3249.202748323613
$ws.Range("R4").Value = 29242.82473491251
$ws.Range("S4").Value = 0.001307161318343143
$ws.Range("T4").Value = 0.001307161318343143
$ws.Range("G5").Value = 25.35940266666667
$ws.Range("H5").Value = 76.078208
$ws.Range("I5").Value = 0.005186643687654987
$ws.Range("J5").Value = 0.005186643687654986
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1667.661682155477
$ws.Range("R5").Value = 15008.9551393993
$ws.Range("S5").Value = 0.000670903914543773
$ws.Range("T5").Value = 0.0006709039145437729
$ws.Range("I6").Value = 0.9837462940761621
$ws.Range("J6").Value = 0.983746294076162
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 700805.4757877274
$ws.Range("R6").Value = 6307249.282089547
$ws.Range("S6").Value = 0.2819355640719596
$ws.Range("T6").Value = 0.2819355640719596
$ws.Range("I7").Value = 0.9837462940761621
$ws.Range("J7").Value = 0.983746294076162
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.3266327807409862
$ws.Range("T7").Value = 0.3266327807409861
$ws.Range("I8").Value = 0.9837462940761621
$ws.Range("J8").Value = 0.983746294076162
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 616273.5199977858
$ws.Range("R8").Value = 5546461.679980072
$ws.Range("S8").Value = 0.2479281747733036
$ws.Range("T8").Value = 0.2479281747733036
$ws.Range("I9").Value = 0.9837462940761621
$ws.Range("J9").Value = 0.983746294076162
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 316303.9719690107
$ws.Range("R9").Value = 2846735.747721097
$ws.Range("S9").Value = 0.1272497744899128
$ws.Range("T9").Value = 0.1272497744899128
$ws.Range("G10").Value = 51.27300266666666
$ws.Range("H10").Value = 153.819008
$ws.Range("I10").Value = 0.01048663484403512
$ws.Range("J10").Value = 0.01048663484403512
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 7470.514669829319
$ws.Range("R10").Value = 67234.63202846388
$ws.Range("S10").Value = 0.003005404267109556
$ws.Range("T10").Value = 0.003005404267109556
$ws.Range("G11").Value = 51.27300266666666
$ws.Range("H11").Value = 153.819008
$ws.Range("I11").Value = 0.01048663484403512
$ws.Range("J11").Value = 0.01048663484403512
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 8654.867604960542
$ws.Range("R11").Value = 77893.80844464488
$ws.Range("S11").Value = 0.003481872023659511
$ws.Range("T11").Value = 0.003481872023659511
$ws.Range("G12").Value = 51.27300266666666
$ws.Range("H12").Value = 153.819008
$ws.Range("I12").Value = 0.01048663484403512
$ws.Range("J12").Value = 0.01048663484403512
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 6569.412669893746
$ws.Range("R12").Value = 59124.71402904371
$ws.Range("S12").Value = 0.002642888976610943
$ws.Range("T12").Value = 0.002642888976610943
$ws.Range("G13").Value = 51.27300266666666
$ws.Range("H13").Value = 153.819008
$ws.Range("I13").Value = 0.01048663484403512
$ws.Range("J13").Value = 0.01048663484403512
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 3371.767978929877
$ws.Range("R13").Value = 30345.91181036889
$ws.Range("S13").Value = 0.001356469576655117
$ws.Range("T13").Value = 0.001356469576655117
$ws.Range("G14").Value = 2.837922333333333
$ws.Range("H14").Value = 8.513767
$ws.Range("I14").Value = 0.0005804273921477663
$ws.Range("J14").Value = 0.0005804273921477662
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 413.4873972728309
$ws.Range("R14").Value = 3721.386575455479
$ws.Range("S14").Value = 0.0001663468774351771
$ws.Range("T14").Value = 0.0001663468774351771
$ws.Range("G15").Value = 2.837922333333333
$ws.Range("H15").Value = 8.513767
$ws.Range("I15").Value = 0.0005804273921477663
$ws.Range("J15").Value = 0.0005804273921477662
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 479.0404460577596
$ws.Range("R15").Value = 4311.364014519836
$ws.Range("S15").Value = 0.0001927190112502582
$ws.Range("T15").Value = 0.0001927190112502582
$ws.Range("G16").Value = 2.837922333333333
$ws.Range("H16").Value = 8.513767
$ws.Range("I16").Value = 0.0005804273921477663
$ws.Range("J16").Value = 0.0005804273921477662
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 363.6120758126543
$ws.Range("R16").Value = 3272.508682313888
$ws.Range("S16").Value = 0.0001462819273527887
$ws.Range("T16").Value = 0.0001462819273527886
$ws.Range("G17").Value = 2.837922333333333
$ws.Range("H17").Value = 8.513767
$ws.Range("I17").Value = 0.0005804273921477663
$ws.Range("J17").Value = 0.0005804273921477662
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 186.6248347581976
$ws.Range("R17").Value = 1679.623512823779
$ws.Range("S17").Value = 0.00007507957610954235
$ws.Range("T17").Value = 0.00007507957610954233
